# Update output for final week
# Refreshes the diamonds table1 summary (Overall / X0..not.D.E.F / X1..best.D.E.F)
# with the latest counts/percentages/means now that more records have come in.
# "Overall" N stays at 150 (81 + 69 = 150, same as 72 + 78 previously); only the
# split between the two groups (and the stats computed within each) changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- X0..not.D.E.F column (C) ---
# n : 72 -> 81 (looks numeric, force text so the leading-space padding survives)
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "     81"

$ws.Range("C4").Value  = "      2 ( 2.5) "
$ws.Range("C5").Value  = "      6 ( 7.4) "
$ws.Range("C6").Value  = "     17 (21.0) "
$ws.Range("C7").Value  = "     14 (17.3) "
$ws.Range("C8").Value  = "     42 (51.9) "
$ws.Range("C9").Value  = "   0.89 (0.59)"
$ws.Range("C10").Value = "4557.68 (4972.55)"

# --- X1..best.D.E.F column (D) ---
# n : 78 -> 69
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "     69"

$ws.Range("D4").Value  = "      2 ( 2.9) "
$ws.Range("D5").Value  = "      9 (13.0) "
$ws.Range("D6").Value  = "     16 (23.2) "
$ws.Range("D7").Value  = "     15 (21.7) "
$ws.Range("D8").Value  = "     27 (39.1) "
$ws.Range("D9").Value  = "   0.75 (0.48)"
$ws.Range("D10").Value = "4120.88 (4785.65)"

# --- Overall column (B) ---
# n stays "    150" (81 + 69 = 150), unchanged
$ws.Range("B4").Value  = "      4 ( 2.7) "
$ws.Range("B5").Value  = "     15 (10.0) "
$ws.Range("B6").Value  = "     33 (22.0) "
$ws.Range("B7").Value  = "     29 (19.3) "
$ws.Range("B8").Value  = "     69 (46.0) "
$ws.Range("B9").Value  = "   0.82 (0.54)"
$ws.Range("B10").Value = "4356.75 (4876.03)"
